$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-like numeric price cells to stay stored as text (matches source data format)
$textCells = @("D8", "D10", "D16", "D17", "D20", "D22", "D25", "D30", "D33", "D36", "D38", "D40", "D41", "D43", "D45", "D46", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.078.78"
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("D3").Value = "1.691.46"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "24.18"
$ws.Range("E8").Value = "  +5.76%  "
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("D10").Value = "0.0627"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D12").Value = "1.931.87"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "1.689.39"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").Value = "66.97"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "250.33"
$ws.Range("E17").Value = "  +6.26%  "
$ws.Range("D18").Value = "28.052.14"
$ws.Range("E18").Value = "  +3.21%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "7.71"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "4.55"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").Value = "147.82"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "3.19"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.452.16"
$ws.Range("E34").Value = "  -5.99%  "
$ws.Range("D36").Value = "0.949"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "0.593"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").Value = "69.38"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "5.58"
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("D44").Value = "1.837.36"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "2.23"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "0.796"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +7.35%  "
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "7.99"
$ws.Range("E51").Value = "  -3.23%  "
